$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Planilha_02"

$ws2.Cells.Item(1, 1).Value = "CAMPO_01"
$ws2.Cells.Item(1, 2).Value = "CAMPO_02"
$ws2.Cells.Item(1, 3).Value = "CAMPO_02"
$ws2.Cells.Item(2, 1).Value = "VALOR_01_01"
$ws2.Cells.Item(2, 2).Value = "VALOR_02_1"
$ws2.Cells.Item(2, 3).Value = "VALOR_03_1"
$ws2.Cells.Item(3, 1).Value = "VALOR_01_2"
$ws2.Cells.Item(3, 2).Value = "VALOR_02_2"
$ws2.Cells.Item(3, 3).Value = "VALOR_03_2"
$ws2.Cells.Item(4, 1).Value = "VALOR_01_3"
$ws2.Cells.Item(4, 2).Value = "VALOR_02_3"
$ws2.Cells.Item(4, 3).Value = "VALOR_03_3"
$ws2.Cells.Item(5, 1).Value = "VALOR_01_4"
$ws2.Cells.Item(5, 2).Value = "VALOR_02_4"
$ws2.Cells.Item(5, 3).Value = "VALOR_03_4"
$ws2.Cells.Item(6, 1).Value = "VALOR_01_5"
$ws2.Cells.Item(6, 2).Value = "VALOR_02_5"
$ws2.Cells.Item(6, 3).Value = "VALOR_03_5"
$ws2.Cells.Item(7, 1).Value = "VALOR_01_6"
$ws2.Cells.Item(7, 2).Value = "VALOR_02_6"
$ws2.Cells.Item(7, 3).Value = "VALOR_03_6"
$ws2.Cells.Item(8, 1).Value = "VALOR_01_7"
$ws2.Cells.Item(8, 2).Value = "VALOR_02_7"
$ws2.Cells.Item(8, 3).Value = "VALOR_03_7"
$ws2.Cells.Item(9, 1).Value = "VALOR_01_8"
$ws2.Cells.Item(9, 2).Value = "VALOR_02_8"
$ws2.Cells.Item(9, 3).Value = "VALOR_03_8"
$ws2.Cells.Item(10, 1).Value = "VALOR_01_9"
$ws2.Cells.Item(10, 2).Value = "VALOR_02_9"
$ws2.Cells.Item(10, 3).Value = "VALOR_03_9"
$ws2.Cells.Item(11, 1).Value = "VALOR_01_10"
$ws2.Cells.Item(11, 2).Value = "VALOR_02_10"
$ws2.Cells.Item(11, 3).Value = "VALOR_03_10"
$ws2.Cells.Item(12, 1).Value = "VALOR_01_11"
$ws2.Cells.Item(12, 2).Value = "VALOR_02_11"
$ws2.Cells.Item(12, 3).Value = "VALOR_03_11"
$ws2.Cells.Item(13, 1).Value = "VALOR_01_12"
$ws2.Cells.Item(13, 2).Value = "VALOR_02_12"
$ws2.Cells.Item(13, 3).Value = "VALOR_03_12"
$ws2.Cells.Item(14, 1).Value = "VALOR_01_13"
$ws2.Cells.Item(14, 2).Value = "VALOR_02_13"
$ws2.Cells.Item(14, 3).Value = "VALOR_03_13"

$ws1.Activate()
